# atKarma Presentation_08272016.pptx — "Update new slide on powerpoint"
#
# Adds a new slide at the end of the deck (after "Issues we encountered"),
# using the same Title+Content layout (slideLayout2 / "obj"), titled
# "API Wishes" with a body paragraph describing the API.ai wish-list item.

$p = $ppt.ActivePresentation

# Same custom layout ("Title and Content") used by the preceding slide.
$newIndex = $p.Slides.Count + 1
$s = $p.Slides.Add($newIndex, 2)

# Title placeholder.
$s.Shapes.Item(1).TextFrame.TextRange.Text = "API Wishes"

# Body / content placeholder — build it as four runs (matching how the
# authoring app split the paragraph around the "API.ai" term).
$tr = $s.Shapes.Item(2).TextFrame.TextRange
$tr.Text = "Users are able to input tasks they want to help others with. "
$tr.InsertAfter("API.ai") | Out-Null
$tr.InsertAfter(" will be able to spot key words and search for specific tasks accordingly based ") | Out-Null
$tr.InsertAfter("on categorization.") | Out-Null
